$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.824.04"
$ws.Range("E2").Value = "  -4.79%  "
$ws.Range("D3").Value = "3.214.03"
$ws.Range("E3").Value = "  -8.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.37%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.205.76"
$ws.Range("E8").Value = "  -8.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -10.52%  "
$ws.Range("E10").Value = "  -10.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.498"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -14.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -15.03%  "
$ws.Range("E14").Value = "  -11.41%  "
$ws.Range("D15").Value = "3.735.58"
$ws.Range("E15").Value = "  -8.24%  "
$ws.Range("D16").Value = "66.863.12"
$ws.Range("E16").Value = "  -4.69%  "
$ws.Range("D17").Value = "3.215.20"
$ws.Range("E17").Value = "  -8.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.115"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "533.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -13.18%  "
$ws.Range("E20").Value = "  -13.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -14.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.763"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -13.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -12.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -10.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -13.48%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -13.47%  "
$ws.Range("E28").Value = "  -14.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -13.75%  "
$ws.Range("E31").Value = "  -10.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "549.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -14.90%  "
$ws.Range("E34").Value = "  -18.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -15.75%  "
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0426"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -10.53%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0866"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -12.89%  "
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -12.42%  "
$ws.Range("E41").Value = "  -12.14%  "
$ws.Range("D42").Value = "2.918.85"
$ws.Range("E42").Value = "  -13.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -23.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.265"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -14.41%  "
$ws.Range("D45").Value = "0.0₃0587"
$ws.Range("E45").Value = "  -20.05%  "
$ws.Range("E46").Value = "  -16.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -16.69%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -16.07%  "
$ws.Range("E50").Value = "  -12.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.72%  "
